# Atualiza a coluna "Salas Preferenciais" (B) com a solucao final para 2023.2
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value  = "206-B,207-B,208-B,209-B,312-C"
$ws.Range("B3").Value  = "201-B,302-B,303-B,304-B,305-B,105-C"
$ws.Range("B4").Value  = "305-B,308-B,309-B,310-B,313-C"
$ws.Range("B5").Value  = "201-A,302-A,303-A,309-A,302-C"
$ws.Range("B6").Value  = "301-A,305-A,307-A,308-A,309-A,105-C"
$ws.Range("B7").Value  = "201-A,202-A,203-A,204-A,205-A,302-C"
$ws.Range("B8").Value  = "301-A,304-A,301-C"
$ws.Range("B9").Value  = "301-B,302-B,303-B,313-C"
$ws.Range("B10").Value = "305-A,306-A,307-A,308-A,313-C"
$ws.Range("B11").Value = "206-A,207-A,208-A,209-A,220-C"
$ws.Range("B12").Value = "201-A,202-A,203-A,204-A,205-A,312-C"
$ws.Range("B13").Value = "303-A,304-A,309-B,310-B,208-B,209-B,301-B,302-B"
$ws.Range("B14").Value = "201-B,202-B,203-B,204-B,205-B,313-C"

[void]$ws.Range("B22").Select()
